$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-09 Friday" "2026-01-10 Saturday"

Replace-Text "785÷6=" "478÷7="
Replace-Text "387÷9=" "186÷7="
Replace-Text "397÷2=" "293÷5="
Replace-Text "143÷3=" "753÷9="
Replace-Text "418÷4=" "927÷2="

Replace-Text "484÷8=" "896÷7="
Replace-Text "956÷5=" "515÷8="
Replace-Text "660÷5=" "976÷5="
Replace-Text "675÷6=" "136÷6="
Replace-Text "396÷3=" "676÷5="

Replace-Text "881÷7=" "669÷4="
Replace-Text "841÷6=" "212÷5="
Replace-Text "424÷9=" "439÷6="
Replace-Text "401÷5=" "510÷6="
Replace-Text "877÷5=" "202÷9="

Replace-Text "807÷8=" "432÷3="
Replace-Text "823÷8=" "228÷8="
Replace-Text "984÷2=" "563÷2="
Replace-Text "410÷4=" "832÷8="
Replace-Text "249÷3=" "879÷9="

Replace-Text "799÷7=" "375÷3="
Replace-Text "407÷7=" "437÷9="
Replace-Text "812÷9=" "143÷6="
Replace-Text "624÷3=" "444÷7="
Replace-Text "174÷6=" "971÷8="
